{"js": "const replacements = [\n  [\"2023-09-21 Thursday\", \"2023-09-22 Friday\"],\n  [\"27\u00d797=\", \"60\u00d762=\"],\n  [\"32\u00d711=\", \"22\u00d730=\"],\n  [\"77\u00d792=\", \"67\u00d774=\"],\n  [\"85\u00d786=\", \"50\u00d758=\"],\n  [\"31\u00d766=\", \"92\u00d752=\"],\n  [\"14\u00d760=\", \"73\u00d798=\"],\n  [\"79\u00d764=\", \"55\u00d773=\"],\n  [\"82\u00d756=\", \"37\u00d793=\"],\n  [\"99\u00d751=\", \"61\u00d768=\"],\n  [\"75\u00d711=\", \"82\u00d754=\"],\n  [\"36\u00d759=\", \"17\u00d731=\"],\n  [\"13\u00d748=\", \"59\u00d751=\"],\n  [\"47\u00d771=\", \"84\u00d755=\"],\n  [\"18\u00d764=\", \"45\u00d794=\"],\n  [\"70\u00d798=\", \"35\u00d783=\"],\n  [\"43\u00d741=\", \"82\u00d737=\"],\n  [\"63\u00d745=\", \"32\u00d778=\"],\n  [\"13\u00d717=\", \"60\u00d798=\"],\n  [\"12\u00d717=\", \"14\u00d767=\"],\n  [\"30\u00d776=\", \"98\u00d743=\"],\n  [\"21\u00d763=\", \"55\u00d797=\"],\n  [\"70\u00d752=\", \"72\u00d743=\"],\n  [\"64\u00d781=\", \"87\u00d713=\"],\n  [\"53\u00d726=\", \"89\u00d770=\"],\n  [\"74\u00d776=\", \"71\u00d793=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Replace each multiplication problem / date string with its updated value.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2023-09-21 Thursday\", \"2023-09-22 Friday\")\n    ,@(\"27\u00d797=\", \"60\u00d762=\")\n    ,@(\"32\u00d711=\", \"22\u00d730=\")\n    ,@(\"77\u00d792=\", \"67\u00d774=\")\n    ,@(\"85\u00d786=\", \"50\u00d758=\")\n    ,@(\"31\u00d766=\", \"92\u00d752=\")\n    ,@(\"14\u00d760=\", \"73\u00d798=\")\n    ,@(\"79\u00d764=\", \"55\u00d773=\")\n    ,@(\"82\u00d756=\", \"37\u00d793=\")\n    ,@(\"99\u00d751=\", \"61\u00d768=\")\n    ,@(\"75\u00d711=\", \"82\u00d754=\")\n    ,@(\"36\u00d759=\", \"17\u00d731=\")\n    ,@(\"13\u00d748=\", \"59\u00d751=\")\n    ,@(\"47\u00d771=\", \"84\u00d755=\")\n    ,@(\"18\u00d764=\", \"45\u00d794=\")\n    ,@(\"70\u00d798=\", \"35\u00d783=\")\n    ,@(\"43\u00d741=\", \"82\u00d737=\")\n    ,@(\"63\u00d745=\", \"32\u00d778=\")\n    ,@(\"13\u00d717=\", \"60\u00d798=\")\n    ,@(\"12\u00d717=\", \"14\u00d767=\")\n    ,@(\"30\u00d776=\", \"98\u00d743=\")\n    ,@(\"21\u00d763=\", \"55\u00d797=\")\n    ,@(\"70\u00d752=\", \"72\u00d743=\")\n    ,@(\"64\u00d781=\", \"87\u00d713=\")\n    ,@(\"53\u00d726=\", \"89\u00d770=\")\n    ,@(\"74\u00d776=\", \"71\u00d793=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
